$wb = $excel.ActiveWorkbook

# --- Sheet "data": add column AA for the new survey wave (30. 3. 2021) ---
$ws1 = $wb.Worksheets.Item("data")
$ws1.Range("AA1").Value = "30. 3. 2021"
$ws1.Range("AA1").Font.Bold = $true
$ws1.Range("AA1").HorizontalAlignment = -4108  # xlCenter
$ws1.Range("AA1").VerticalAlignment = -4160    # xlTop
$ws1.Range("AA1").Borders.LineStyle = 1        # xlContinuous, matches header style of Z1

$ws1.Range("AA2").Value = 0.48
$ws1.Range("AA3").Value = 0.32
$ws1.Range("AA4").Value = 0.2
$ws1.Range("AA5").Value = 0.25
$ws1.Range("AA6").Value = 0.33
$ws1.Range("AA7").Value = 0.42
$ws1.Range("AA8").Value = 0.52
$ws1.Range("AA9").Value = 0.31
$ws1.Range("AA10").Value = 0.17
$ws1.Range("AA11").Value = 0.51
$ws1.Range("AA12").Value = 0.32
$ws1.Range("AA13").Value = 0.17
$ws1.Range("AA14").Value = 0.38
$ws1.Range("AA15").Value = 0.3
$ws1.Range("AA16").Value = 0.32
$ws1.Range("AA17").Value = 0.52
$ws1.Range("AA18").Value = 0.32
$ws1.Range("AA19").Value = 0.16
$ws1.Range("AA20").Value = 0.4
$ws1.Range("AA21").Value = 0.3
$ws1.Range("AA22").Value = 0.3
$ws1.Range("AA23").Value = 0.38
$ws1.Range("AA24").Value = 0.34
$ws1.Range("AA25").Value = 0.28
$ws1.Range("AA26").Value = 0.39
$ws1.Range("AA27").Value = 0.39
$ws1.Range("AA28").Value = 0.22
$ws1.Range("AA29").Value = 0.48
$ws1.Range("AA30").Value = 0.28
$ws1.Range("AA31").Value = 0.24
$ws1.Range("AA32").Value = 0.62
$ws1.Range("AA33").Value = 0.26
$ws1.Range("AA34").Value = 0.12
$ws1.Range("AA35").Value = 0.35
$ws1.Range("AA36").Value = 0.38
$ws1.Range("AA37").Value = 0.27
$ws1.Range("AA38").Value = 0.38
$ws1.Range("AA39").Value = 0.39
$ws1.Range("AA40").Value = 0.23
$ws1.Range("AA41").Value = 0.57
$ws1.Range("AA42").Value = 0.24
$ws1.Range("AA43").Value = 0.19
$ws1.Range("AA44").Value = 0.59
$ws1.Range("AA45").Value = 0.28
$ws1.Range("AA46").Value = 0.13
$ws1.Range("AA47").Value = 0.41
$ws1.Range("AA48").Value = 0.31
$ws1.Range("AA49").Value = 0.28
$ws1.Range("AA50").Value = 0.67
$ws1.Range("AA51").Value = 0.27
$ws1.Range("AA52").Value = 0.06
$ws1.Range("AA53").Value = 0.42
$ws1.Range("AA54").Value = 0.38
$ws1.Range("AA55").Value = 0.2
$ws1.Range("AA56").Value = 0.67
$ws1.Range("AA57").Value = 0.19
$ws1.Range("AA58").Value = 0.14
$ws1.Range("AA59").Value = 0.65
$ws1.Range("AA60").Value = 0.2
$ws1.Range("AA61").Value = 0.15

# --- Sheet "pocetR": add column Z for the new survey wave (30. 3. 2021) ---
$ws2 = $wb.Worksheets.Item("pocetR")
$ws2.Range("Z1").Value = "30. 3. 2021"
$ws2.Range("Z1").Font.Bold = $true
$ws2.Range("Z1").HorizontalAlignment = -4108  # xlCenter
$ws2.Range("Z1").VerticalAlignment = -4160    # xlTop
$ws2.Range("Z1").Borders.LineStyle = 1        # xlContinuous, matches header style of Y1

$ws2.Range("Z2").Value = 1142
$ws2.Range("Z3").Value = 146
$ws2.Range("Z4").Value = 996
$ws2.Range("Z5").Value = 898
$ws2.Range("Z6").Value = 175
$ws2.Range("Z7").Value = 9
$ws2.Range("Z8").Value = 60
$ws2.Range("Z9").Value = 859
$ws2.Range("Z10").Value = 161
$ws2.Range("Z11").Value = 71
$ws2.Range("Z12").Value = 51
$ws2.Range("Z13").Value = 410
$ws2.Range("Z14").Value = 456
$ws2.Range("Z15").Value = 276
$ws2.Range("Z16").Value = 132
$ws2.Range("Z17").Value = 324
$ws2.Range("Z18").Value = 385
$ws2.Range("Z19").Value = 181
$ws2.Range("Z20").Value = 316
$ws2.Range("Z21").Value = 95
$ws2.Range("Z22").Value = 304
$ws2.Range("Z23").Value = 168
$ws2.Range("Z24").Value = 98
$ws2.Range("Z25").Value = ""

# --- Update "aktualizace" date in the two summary description strings ---
$ws1.Range("A62").Value = "Život během pandemie, Obavy ze ztráty práce, % respondentů celkově a ve skupinách, aktualizace 7. 4. 2021"
$ws2.Range("A25").Value = "Život během pandemie, Obavy ze ztráty práce, velikost dotázaného souboru celkově a ve skupinách, aktualizace 7. 4. 2021"

